# implemented analog input signal type
#
# Adds a new "ADC_100MS" worksheet (after the existing "MEM1" sheet)
# describing the analog-input pin table (supplies/grounds + two analog
# input pins), and moves the active sheet/selection over to it.

$wb  = $excel.ActiveWorkbook
$mem1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Create the new sheet right after MEM1.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $mem1)
$ws.Name = "ADC_100MS"

# Column widths (mirrors MEM1's layout: wider pin-name / type columns,
# a wide comment column, narrow everything else).
$ws.Columns.Item(4).ColumnWidth = 15.636666666666665   # D ~16.47
$ws.Columns.Item(5).ColumnWidth = 16.526666666666667   # E ~17.36
$ws.Columns.Item(8).ColumnWidth = 25.55666666666667    # H ~26.39

# Page setup / margins - copy MEM1's print layout onto the new sheet.
$ps = $ws.PageSetup
$ps.LeftMargin = 56.7
$ps.RightMargin = 56.7
$ps.TopMargin = 75.8
$ps.BottomMargin = 75.8
$ps.HeaderMargin = 56.7
$ps.FooterMargin = 56.7
$ps.PaperSize = 9
$ps.Zoom = 100
$ps.Orientation = 1
$ps.CenterHeader = "&""Times New Roman,Regular""&12&A"
$ps.CenterFooter = "&""Times New Roman,Regular""&12Page &P"

# A1 carries only centered formatting (no value) - matches the
# source sheet's stray formatted cell.
$ws.Range("A1").HorizontalAlignment = -4108            # xlCenter

# ---------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "pin"
$ws.Range("C5").Value = "type"
$ws.Range("D5").Value = "related ground"
$ws.Range("E5").Value = "related supply"
$ws.Range("F5").Value = "x"
$ws.Range("G5").Value = "y"
$ws.Range("H5").Value = "comment"
$ws.Range("F5").HorizontalAlignment = -4108
$ws.Range("G5").HorizontalAlignment = -4108

# Blank centered spacer cells under the x/y header (no values).
$ws.Range("F6").HorizontalAlignment = -4108
$ws.Range("G6").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3. Analog / digital supply + ground rows.
# ---------------------------------------------------------------------
$ws.Range("B7").Value = "vdda"
$ws.Range("C7").Value = "supply"
$ws.Range("D7").Value = "n.a."
$ws.Range("E7").Value = "n.a."
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 10
$ws.Range("H7").Value = "analog supply"
$ws.Range("F7").HorizontalAlignment = -4108
$ws.Range("G7").HorizontalAlignment = -4108

$ws.Range("B8").Value = "gnda"
$ws.Range("C8").Value = "ground"
$ws.Range("D8").Value = "n.a."
$ws.Range("E8").Value = "n.a."
$ws.Range("F8").Formula = "=F7"
$ws.Range("G8").Formula = "=G7+10"
$ws.Range("H8").Value = "analog ground"
$ws.Range("F8").HorizontalAlignment = -4108
$ws.Range("G8").HorizontalAlignment = -4108

$ws.Range("B9").Value = "vddd"
$ws.Range("C9").Value = "supply"
$ws.Range("D9").Value = "n.a."
$ws.Range("E9").Value = "n.a."
$ws.Range("F9").Formula = "=F8"
$ws.Range("G9").Formula = "=G8+10"
$ws.Range("H9").Value = "digital supply"
$ws.Range("F9").HorizontalAlignment = -4108
$ws.Range("G9").HorizontalAlignment = -4108

$ws.Range("B10").Value = "gndd"
$ws.Range("C10").Value = "ground"
$ws.Range("D10").Value = "n.a."
$ws.Range("E10").Value = "n.a."
$ws.Range("F10").Formula = "=F9"
$ws.Range("G10").Formula = "=G9+10"
$ws.Range("H10").Value = "digital ground"
$ws.Range("F10").HorizontalAlignment = -4108
$ws.Range("G10").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 4. Analog input pins (row 11 stays blank, matching MEM1's spacing).
# ---------------------------------------------------------------------
$ws.Range("B12").Value = "v_in1"
$ws.Range("C12").Value = "ana_in"
$ws.Range("D12").Value = "gnda"
$ws.Range("E12").Value = "vdda"
$ws.Range("F12").Value = 20
$ws.Range("G12").Value = 10
$ws.Range("H12").Value = "input voltage 1"
$ws.Range("F12").HorizontalAlignment = -4108
$ws.Range("G12").HorizontalAlignment = -4108

$ws.Range("B13").Value = "v_in2"
$ws.Range("C13").Value = "ana_in"
$ws.Range("D13").Value = "gnda"
$ws.Range("E13").Value = "vdda"
$ws.Range("F13").Formula = "=F12+10"
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = "input voltage 2"
$ws.Range("F13").HorizontalAlignment = -4108
$ws.Range("G13").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 5. Selection bookkeeping: MEM1 loses the tab focus / its selection
#    resets to B6; ADC_100MS becomes active with G18 selected.
# ---------------------------------------------------------------------
$mem1.Range("B6").Select()
$ws.Activate()
$ws.Range("G18").Select()
